$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = "No code"

[void]$ws.Range("A3").Select()
